$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124; existing rows 124-133 shift down to 125-134.
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new record (Lemon / Segunda).
$ws.Range("A124").Value = 4
$ws.Range("B124").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C124").Value = "Los Lagos"
$ws.Range("D124").Value = 44585
$ws.Range("E124").Value = 10
$ws.Range("F124").Value = "Fruta"
$ws.Range("G124").Value = 100103
$ws.Range("H124").Value = "Frutos de hueso (carozo)"
$ws.Range("I124").Value = 100103002
$ws.Range("J124").Value = "Ciruela"
$ws.Range("K124").Value = "Lemon"
$ws.Range("L124").Value = "Segunda"
$ws.Range("M124").Value = 200
$ws.Range("N124").Value = 14000
$ws.Range("O124").Value = 14000
$ws.Range("P124").Value = 14000
$ws.Range("Q124").Value = "$/caja 15 kilos granel"
$ws.Range("R124").Value = "Región de O'Higgins"
$ws.Range("S124").Value = 933
$ws.Range("T124").Value = 15
